# Apply the commit: re-run the Neo4j stats query for breed
# 'Flat-Coated Retriever' instead of 'Akita', and zero out the
# resulting file/sample counts on the StatOutput sheet.

$wb = $excel.ActiveWorkbook

$statOutput = $wb.Worksheets.Item("StatOutput")
$statMessage = $wb.Worksheets.Item("StatOutput_Message")

# New Cypher query text (same as before but filtering on
# 'Flat-Coated Retriever' rather than 'Akita').
$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Flat-Coated Retriever']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Update the Cypher query recorded in the message log (row 18).
$statMessage.Range("A18").Value = $newCypher

# Update the resulting counts on the StatOutput sheet: the new
# breed filter returns zero files and zero samples. These are
# stored as text (shared strings) elsewhere on the sheet, so
# force text via a leading apostrophe rather than changing the
# cell's number format.
$statOutput.Range("A2").Value = "'0"
$statOutput.Range("B2").Value = "'0"
